$wb = $excel.ActiveWorkbook

# --- Update the conversion note text on "Hoja1" ---
$ws1 = $wb.Worksheets("Hoja1")
$ws1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 7.05 = 28379.14 pesos`n✅ 28379.14 pesos = 7.04 = 953.48 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

# --- Update the rate values on "tasas" sheet ---
$ws2 = $wb.Worksheets("tasas")
$ws2.Range("N10").Value = 141.9
$ws2.Range("O10").Value = 4027
$ws2.Range("N12").Value = 4030
$ws2.Range("O12").Value = 135.4
